$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCtMbCtbDP")
$ws.Activate()

# B1 used to hold the numeric calibration value (0.8); it now holds the
# row's label text instead.
$ws.Range("B1").Value = "share of costs that must be covered"

# The per-technology calibration values in B2:B24 move from 0.8 to 1.
$ws.Range("B2:B24").Value = 1

# Reflect the new selection left on the sheet (B2:B24, active cell B2).
$ws.Range("B2:B24").Select()
